$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new study ("68 - Nigro et al. 2019") was added to the extraction sheet.
# It belongs alphabetically/numerically between row 39 (Nakanishi et al. 1991)
# and the former row 40 (Pollok et al. 2009), so insert a fresh row at 40 and
# push everything below it down by one.
$ws.Rows(40).Insert()

# Column A ("68") looks numeric, so force the cell to Text first - otherwise
# Excel would store it as a number instead of the shared-string label used by
# every other row in this column.
$ws.Cells.Item(40, 1).NumberFormat = "@"
$ws.Cells.Item(40, 1).Value = "68"
$ws.Cells.Item(40, 1).Style = "Normal"

$ws.Cells.Item(40, 2).Value = "Nigro et al. 2019"
$ws.Cells.Item(40, 3).Value = "Quasi-experimental study"
$ws.Cells.Item(40, 4).Value = "UPDRS_III_item20__apomorphine_baseline"
$ws.Cells.Item(40, 5).Value = "UPDRS_III_item20__apomorphine_follow-up"
$ws.Cells.Item(40, 6).Value = "UPDRS_III_item20__placebo_baseline"
$ws.Cells.Item(40, 7).Value = "UPDRS_III_item20__placebo_follow-up"
